$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$objetivosText = @'
Demonstrar as principais etapas no desenvolvimento dos processos bioquímicos industriais abordando aspectos bioquímicos importantes na produção de alimentos, e importantes metabólitos. Apresentar aos alunos uma visão das aplicações potenciais e estratégicas da biotecnologia moderna, incluindo as biorefinarias de lignocelulósicos.Aprimorar o raciocínio e despertar o espírito crítico e a criatividade dos alunos na resolução de problemas industriais envolvendo processos bioquímicos.
'@
$programaResumidoText = @'
Introdução aos processos bioquímicos industriais que incluem o processamento de alimentos, e importantes metabólitos, a manufatura de soros e vacinas, e os conceitos modernos de bioenergia e biorrefinarias.
'@
$programaText = @'
1. Introdução ao Processamento de alimentos: tipos de indústria de alimentos, matériasprimas,fases doprocessamento de produtos alimentícios, conservação/alterações de alimentos, microbiologia dealimentos, alterações bioquímicas em alimentos (oxidação de lipídeos, antioxidantes, escurecimentoenzimático e não enzimático), aflatoxinas, conservantes químicos, toxicantes naturais.2. Discussão e apresentação sobre aspectos bioquímicos importantes na produção de metabólitos por microrganismos e  estudo de casos .3. Manufatura de soros e vacinas Métodosindustriais para a produção de soros e vacinas 4.Biotecnologia de lignocelulósicos: Separação e fermentação das frações e principais processosbioquímicos envolvendo materiais lignocelulósicos.5. Bioenergia, biocombustíveis e biorrefinarias.
'@
$metodoText = @'
A avaliação será feita por meio de prova escrita (P1) e trabalhos (P2).
'@
$criterioText = @'
A nota final (NF) será calculada da seguintes maneira: NF=(P1 + P2)/2
'@
$normaRecuperacaoText = @'
A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR)será calculada como MR=(NF+PR)/2
'@
$bibliografiaText = @'
1. BORZANI, W., SCHMIDELL, W., LIMA, U.A., AQUARONE, E. Biotecnologia Industrial - Fundamentos (Vol 1). São Paulo: Edgard Blucher Ltda, 2001.
2. DEMAIN, A.L., SOLOMON, N.A. (Eds). Manual of industrial microbiology and biotechnology, Washington: American Society for Microbiology, 1986.
3. WANG, D.C. et al. Fermentation and Enzyme Technology, New York: Wiley-Interscience, 1979.
4. GAVA, A.J. Princípios de Tecnologia de Alimentos, São Paulo: Nobel, 1983.
5. LIMA , U. A et al. Biotecnología Industrial, Biotecnologia na produção de alimentos - Série Biotecnología, vol4. Ed. Edgard Blucher,Ltda , 2001.
6. EVANGELISTA, J. Tecnologia de Alimentos, Rio de Janeiro: Livraria Atheneu, 1987.
7. CAMARGO R. et al., Tecnologia de produtos Agropecuários- Alimentos, São Paulo: Livraria Nobel, 1984.
'@
$docentesText = @'
1814052 - Silvio Silverio da Silva
'@

# --- 1. Update "Objetivos:" row (row 10) with the new long objective text ---
$ws.Range("B10").Value = $objetivosText
$ws.Range("C10").Value = $objetivosText

# --- 2. Insert a new row at 13 for "Docentes responsaveis:" data (shifts old rows 13-24 down to 14-25) ---
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()
$ws.Range("B13").Value = $docentesText
$ws.Range("C13").Value = $docentesText
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)

# --- 3. Update "Programa resumido:" row (now row 14) ---
$ws.Range("B14").Value = $programaResumidoText
$ws.Range("C14").Value = $programaResumidoText

# --- 4. Update "Programa:" row (now row 16) ---
$ws.Range("B16").Value = $programaText
$ws.Range("C16").Value = $programaText

# --- 5. Update "Metodo:" row (now row 19) ---
$ws.Range("B19").Value = $metodoText
$ws.Range("C19").Value = $metodoText

# --- 6. Update "Criterio:" row (now row 20) ---
$ws.Range("B20").Value = $criterioText
$ws.Range("C20").Value = $criterioText

# --- 7. Update "Norma de recuperacao:" row (now row 21) ---
$ws.Range("B21").Value = $normaRecuperacaoText
$ws.Range("C21").Value = $normaRecuperacaoText

# --- 8. Update "Bibliografia:" row (now row 22) with the new reference list ---
$ws.Range("B22").Value = $bibliografiaText
$ws.Range("C22").Value = $bibliografiaText
